# Daily attendance processing - 2026-01-19 21:02:59
#
# For every data row, column G ("Recorded By") holds a comma-separated list
# of the users/processes that recorded/touched that attendance session
# (e.g. "System, dnasr281@gmail.com"). This pass re-orders the "System"
# token within that list:
#   - if "System" is the first entry, move it to the end of the list
#   - otherwise (it appears later in the list), move it to the front
# All other entries keep their relative order. Cells that don't contain a
# "System" token (or don't contain a comma-separated list at all) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G is "Recorded By" (see header row 1). Work out how many rows are
# in use so we cover the whole table regardless of its exact extent.
$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Rows.Count + $firstRow - 1

$recordedByCol = 7  # column G

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    # Split the "Recorded By" list on the ", " separator used in this sheet.
    $parts = @($val -split ', ')
    if ($parts.Length -lt 2) {
        # Nothing to reorder for a single-entry cell.
        continue
    }

    # Locate the "System" token (case-sensitive; there's also a distinct,
    # deliberately-lowercase "system" entry in a few rows that must NOT be
    # treated the same as "System").
    $idx = -1
    for ($k = 0; $k -lt $parts.Length; $k++) {
        if ($parts[$k].Equals("System")) {
            $idx = $k
        }
    }

    if ($idx -eq -1) {
        # No "System" token present - leave the cell alone.
        continue
    }

    # Rebuild the list without the "System" token...
    $newParts = New-Object System.Collections.ArrayList
    for ($k = 0; $k -lt $parts.Length; $k++) {
        if ($k -ne $idx) {
            [void]$newParts.Add($parts[$k])
        }
    }

    # ...then reinsert it at the opposite end from where it started.
    if ($idx -eq 0) {
        [void]$newParts.Add("System")
    } else {
        [void]$newParts.Insert(0, "System")
    }

    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
        $changed++
    }
}

Write-Output "Reordered 'System' token in $changed 'Recorded By' cells"
